$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2sxc Scripts" snippet: insert a whole row above row 250 (this
# shifts all existing snippet rows below it down by one, same as using the
# Excel UI's "Insert Sheet Rows").
$ws.Rows("250:250").Insert()

# The worksheet's table ("Table1") needs to grow by one row to keep covering
# the newly inserted row.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:G267"))

# Fill in the new snippet row. Column A ("set") uses the same leading
# apostrophe the rest of the "set" column uses, so Excel marks it with the
# quote-prefix cell format like its neighbours.
$ws.Range("A250").Value = "'@Html"
$ws.Range("B250").Value = "2sxc Scripts"
$ws.Range("C250").Value = "standard 2sxc JS"
$ws.Range("E250").Value = '<script type="text/javascript" src="/desktopmodules/tosic_sexycontent/js/2sxc.api.min.js" data-enableoptimizations="100"></script>'
$ws.Range("F250").Value = "a css-tag which enables optimizations - with def. priority at page bottom"
$ws.Range("G250").Value = "read api-docs:https://github.com/2sic/2sxc/wiki/Template-Assets"

# Leave the cursor on the new row's content cell, matching the author's
# recorded selection after making this edit.
$ws.Range("E250").Select()
